$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.886.75"
$ws.Range("E2").Value = "'  -1.22%  "
$ws.Range("D3").Value = "'2.456.54"
$ws.Range("E3").Value = "'  -2.59%  "
$ws.Range("E5").Value = "'  -1.39%  "
$ws.Range("D6").Value = "'168.73"
$ws.Range("E6").Value = "'  -1.86%  "
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("D8").Value = "'0.510"
$ws.Range("E8").Value = "'  -2.33%  "
$ws.Range("D9").Value = "'2.457.32"
$ws.Range("E9").Value = "'  -2.37%  "
$ws.Range("E10").Value = "'  -2.05%  "
$ws.Range("E11").Value = "'  -1.25%  "
$ws.Range("E12").Value = "'  -2.56%  "
$ws.Range("E13").Value = "'  -4.46%  "
$ws.Range("D15").Value = "'25.11"
$ws.Range("E15").Value = "'  -4.22%  "
$ws.Range("D16").Value = "'66.758.51"
$ws.Range("E16").Value = "'  -1.82%  "
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("E17").Value = "'  -3.99%  "
$ws.Range("D18").Value = "'2.465.91"
$ws.Range("E18").Value = "'  -1.88%  "
$ws.Range("D19").Value = "'10.86"
$ws.Range("E19").Value = "'  -8.24%  "
$ws.Range("D20").Value = "'7.33"
$ws.Range("E20").Value = "'  -8.13%  "
$ws.Range("D21").Value = "'348.14"
$ws.Range("E21").Value = "'  -4.43%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "'  -3.28%  "
$ws.Range("E23").Value = "'  +0.25%  "
$ws.Range("D24").Value = "'68.52"
$ws.Range("E24").Value = "'  -4.76%  "
$ws.Range("D25").Value = "'4.17"
$ws.Range("E25").Value = "'  -7.24%  "
$ws.Range("D26").Value = "'1.78"
$ws.Range("E26").Value = "'  -4.69%  "
$ws.Range("D27").Value = "'9.05"
$ws.Range("E27").Value = "'  -7.72%  "
$ws.Range("E28").Value = "'  -38.32%  "
$ws.Range("E29").Value = "'  -3.05%  "
$ws.Range("D30").Value = "'508.95"
$ws.Range("E30").Value = "'  -3.93%  "
$ws.Range("D31").Value = "'0.0₃0889"
$ws.Range("E31").Value = "'  -6.09%  "
$ws.Range("D32").Value = "'7.55"
$ws.Range("E32").Value = "'  -7.98%  "
$ws.Range("E33").Value = "'  -6.09%  "
$ws.Range("D34").Value = "'1.21"
$ws.Range("E34").Value = "'  -6.06%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "'  -0.16%  "
$ws.Range("D36").Value = "'158.17"
$ws.Range("E36").Value = "'  -0.90%  "
$ws.Range("E37").Value = "'  -11.73%  "
$ws.Range("D38").Value = "'18.63"
$ws.Range("E38").Value = "'  +0.21%  "
$ws.Range("D39").Value = "'18.11"
$ws.Range("E39").Value = "'  -5.65%  "
$ws.Range("D40").Value = "'1.31"
$ws.Range("E40").Value = "'  -7.69%  "
$ws.Range("E41").Value = "'  +0.13%  "
$ws.Range("E42").Value = "'  -6.14%  "
$ws.Range("D43").Value = "'4.73"
$ws.Range("E43").Value = "'  -6.28%  "
$ws.Range("E44").Value = "'  -6.34%  "
$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = "'  -3.94%  "
$ws.Range("D46").Value = "'38.50"
$ws.Range("E46").Value = "'  -1.70%  "
$ws.Range("D47").Value = "'140.46"
$ws.Range("E47").Value = "'  -4.76%  "
$ws.Range("D48").Value = "'3.41"
$ws.Range("E48").Value = "'  -7.22%  "
$ws.Range("D49").Value = "'0.505"
$ws.Range("E49").Value = "'  -8.09%  "
$ws.Range("D50").Value = "'0.0726"
$ws.Range("E50").Value = "'  -2.11%  "
$ws.Range("E51").Value = "'  -9.16%  "

Write-Output "Updated cryptos list"
